$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 914.80646
$ws.Range("J17").Value = 914.80646
$ws.Range("L17").Value = 2744.41938
$ws.Range("N17").Value = -3080.41938
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H58").Value = 2177.625
$ws.Range("I58").Value = 333.33334
$ws.Range("J58").Value = 2603.2307
$ws.Range("K58").Value = 1000.00002
$ws.Range("L58").Value = 7809.6921
$ws.Range("M58").Value = -850.0000200000001
$ws.Range("N58").Value = -8109.6921
$ws.Range("H127").Value = 1601.2142
$ws.Range("I127").Value = 800
$ws.Range("J127").Value = 1734.75
$ws.Range("K127").Value = 2400
$ws.Range("L127").Value = 5204.25
$ws.Range("M127").Value = 2560
$ws.Range("N127").Value = -15124.25
$ws.Range("H138").Value = 1378.0509
$ws.Range("I138").Value = 1003.587
$ws.Range("J138").Value = 2703.077
$ws.Range("K138").Value = 3010.761
$ws.Range("L138").Value = 8109.231000000001
$ws.Range("M138").Value = 2129.239
$ws.Range("N138").Value = -18389.231
$ws.Range("M20").Value = ""
$ws.Range("M35").Value = ""

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6620.5
$ws.Range("I32").Value = 6296.698
$ws.Range("J32").Value = 7940.615
$ws.Range("K32").Value = 6296.698
$ws.Range("L32").Value = 7940.615
$ws.Range("M32").Value = -6009.698
$ws.Range("N32").Value = -8514.615
$ws.Range("H37").Value = 5350
$ws.Range("I37").Value = 5350
$ws.Range("K37").Value = 5350
$ws.Range("M37").Value = -5077
$ws.Range("H61").Value = 2838.96
$ws.Range("I61").Value = 1600.9
$ws.Range("J61").Value = 7791.2
$ws.Range("K61").Value = 1600.9
$ws.Range("L61").Value = 7791.2
$ws.Range("M61").Value = -1388.9
$ws.Range("N61").Value = -8215.200000000001
$ws.Range("H63").Value = 4693.5415
$ws.Range("I63").Value = 2664.3333
$ws.Range("J63").Value = 8075.5557
$ws.Range("K63").Value = 2664.3333
$ws.Range("L63").Value = 8075.5557
$ws.Range("M63").Value = -1978.3333
$ws.Range("N63").Value = -9447.555700000001
$ws.Range("H66").Value = 4693.5415
$ws.Range("I66").Value = 2664.3333
$ws.Range("J66").Value = 8075.5557
$ws.Range("K66").Value = 13321.6665
$ws.Range("L66").Value = 40377.7785
$ws.Range("M66").Value = -9889.666499999999
$ws.Range("N66").Value = -47241.7785
$ws.Range("H132").Value = 7743.7856
$ws.Range("I132").Value = 5044.625
$ws.Range("J132").Value = 11342.667
$ws.Range("K132").Value = 15133.875
$ws.Range("L132").Value = 34028.001
$ws.Range("M132").Value = -12603.875
$ws.Range("N132").Value = -39088.001
$ws.Range("H136").Value = 2838.96
$ws.Range("I136").Value = 1600.9
$ws.Range("J136").Value = 7791.2
$ws.Range("K136").Value = 4802.700000000001
$ws.Range("L136").Value = 23373.6
$ws.Range("M136").Value = -2252.700000000001
$ws.Range("N136").Value = -28473.6

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 24090
$ws.Range("J74").Value = 24090
$ws.Range("L74").Value = 24090
$ws.Range("N74").Value = -25962
$ws.Range("H77").Value = 24090
$ws.Range("J77").Value = 24090
$ws.Range("L77").Value = 72270
$ws.Range("N77").Value = -81630
$ws.Range("H134").Value = 5326.2793
$ws.Range("I134").Value = 2505.8096
$ws.Range("K134").Value = 7517.4288
$ws.Range("M134").Value = -4982.4288
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10755964
$ws.Range("I31").Value = 2096.2
$ws.Range("J31").Value = 20837714
$ws.Range("K31").Value = 2096.2
$ws.Range("L31").Value = 20837714
$ws.Range("M31").Value = -1801.2
$ws.Range("N31").Value = -20838304
$ws.Range("H34").Value = 10755964
$ws.Range("I34").Value = 2096.2
$ws.Range("J34").Value = 20837714
$ws.Range("K34").Value = 2096.2
$ws.Range("L34").Value = 20837714
$ws.Range("M34").Value = -1894.2
$ws.Range("N34").Value = -20838118
$ws.Range("H58").Value = 3212.3906
$ws.Range("I58").Value = 1512.4375
$ws.Range("J58").Value = 8312.25
$ws.Range("K58").Value = 1512.4375
$ws.Range("L58").Value = 8312.25
$ws.Range("M58").Value = -1309.4375
$ws.Range("N58").Value = -8718.25
$ws.Range("H99").Value = 2874.75
$ws.Range("I99").Value = 2735
$ws.Range("J99").Value = 3666.6667
$ws.Range("K99").Value = 2735
$ws.Range("L99").Value = 3666.6667
$ws.Range("M99").Value = -1237
$ws.Range("N99").Value = -6662.6667
$ws.Range("H126").Value = 2874.75
$ws.Range("I126").Value = 2735
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 8205
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -5735
$ws.Range("N126").Value = -15940.0001
$ws.Range("H136").Value = 3212.3906
$ws.Range("I136").Value = 1512.4375
$ws.Range("J136").Value = 8312.25
$ws.Range("K136").Value = 4537.3125
$ws.Range("L136").Value = 24936.75
$ws.Range("M136").Value = -1987.3125
$ws.Range("N136").Value = -30036.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2433.6667
$ws.Range("I35").Value = 300
$ws.Range("J35").Value = 2860.4
$ws.Range("K35").Value = 900
$ws.Range("L35").Value = 8581.200000000001
$ws.Range("M35").Value = -612
$ws.Range("N35").Value = -9157.200000000001
$ws.Range("H55").Value = 3145.5557
$ws.Range("J55").Value = 3416.25
$ws.Range("L55").Value = 10248.75
$ws.Range("N55").Value = -10602.75
$ws.Range("H59").Value = 2500
$ws.Range("J59").Value = 3000
$ws.Range("L59").Value = 9000
$ws.Range("N59").Value = -10080
$ws.Range("H80").Value = 2449.3333
$ws.Range("J80").Value = 2948.9092
$ws.Range("L80").Value = 8846.7276
$ws.Range("N80").Value = -10718.7276
$ws.Range("H83").Value = 2449.3333
$ws.Range("J83").Value = 2948.9092
$ws.Range("L83").Value = 26540.1828
$ws.Range("N83").Value = -35900.1828
$ws.Range("H106").Value = 3849.0908
$ws.Range("J106").Value = 3849.0908
$ws.Range("L106").Value = 11547.2724
$ws.Range("N106").Value = -13439.2724

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 746150.6
$ws.Range("I132").Value = 1192055.5
$ws.Range("J132").Value = 2975.8096
$ws.Range("K132").Value = 3576166.5
$ws.Range("L132").Value = 8927.4288
$ws.Range("M132").Value = -3573636.5
$ws.Range("N132").Value = -13987.4288

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 70005
$ws.Range("J3").Value = 70005
$ws.Range("L3").Value = 70005
$ws.Range("N3").Value = -70229
$ws.Range("H13").Value = 35878
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 47670.668
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 47670.668
$ws.Range("M13").Value = -360
$ws.Range("N13").Value = -47950.668
$ws.Range("H15").Value = 70005
$ws.Range("J15").Value = 70005
$ws.Range("L15").Value = 70005
$ws.Range("N15").Value = -70345
$ws.Range("H16").Value = 2883.75
$ws.Range("I16").Value = 2010
$ws.Range("K16").Value = 2010
$ws.Range("M16").Value = -1840
$ws.Range("H132").Value = 26345306
$ws.Range("I132").Value = 28602990
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 85808970
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -85806440
$ws.Range("N132").Value = -22058
